# Apply team-roster updates to Sheet1: refresh row 12 (B/C) and append rows 13-34
# with newly scraped club records, extending the used range to A1:C34.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Cells.Item(12, 2).Value = '강원고성군거진FC'
$ws.Cells.Item(12, 3).Value = '창단일자
20210101
주소
강원도 고성군 거진읍 거진운동장길 30 (거진리) 1
연락처
010-4677-7077
감독
팀 관리에서 팀 SNS를 등록하세요'

# Row 13
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(13, 1).Value = 0
$ws.Cells.Item(13, 2).Value = '강원고성군고성사랑축구회'
$ws.Cells.Item(13, 3).Value = '창단일자
20180102
주소
강원 고성군 토성면 장새미1길 11 (용암리) 용암리
연락처
010-2828-4338
감독
팀 관리에서 팀 SNS를 등록하세요'

# Row 14
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(14, 1).Value = 0
$ws.Cells.Item(14, 2).Value = '강원고성군달홀FC'
$ws.Cells.Item(14, 3).Value = '창단일자
20110101
주소
강원 고성군 간성읍 간성로 64 정일빌딩 3층
연락처
010-7708-4342
감독
최서희
팀 관리에서 팀 SNS를 등록하세요'

# Row 15
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15, 1).Value = 0
$ws.Cells.Item(15, 2).Value = '강원고성군죽왕축구회'
$ws.Cells.Item(15, 3).Value = '창단일자
20050101
주소
강원도 고성군 간성읍 탑동길 12 (신안리, 고성경찰서) 금강파출소
연락처
010-9300-5682
감독
강완식
팀 관리에서 팀 SNS를 등록하세요'

# Row 16
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(16, 1).Value = 0
$ws.Cells.Item(16, 2).Value = '강원고성군토요축구회'
$ws.Cells.Item(16, 3).Value = '창단일자
1974
주소
강원 고성군 거진읍 화포2길 26 (화포리)
연락처
010-8727-2188
감독
고경태
팀 관리에서 팀 SNS를 등록하세요'

# Row 17
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17, 1).Value = 0
$ws.Cells.Item(17, 2).Value = '강원동해시초록생활축구회'
$ws.Cells.Item(17, 3).Value = '창단일자
19780120
주소
강원 동해시 중앙로 143 (천곡동, 행복드림아파트) 105동 803호
연락처
010-5574-0690
감독
심영철
팀 관리에서 팀 SNS를 등록하세요'

# Row 18
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(18, 1).Value = 0
$ws.Cells.Item(18, 2).Value = '강원동해시해일생활축구회'
$ws.Cells.Item(18, 3).Value = '창단일자
19990101
주소
강원도 동해시 효자남길 26-1 (동회동, 동보상가) 꿀꿀이뱃살
연락처
010-9728-0717
감독
배기표
팀 관리에서 팀 SNS를 등록하세요'

# Row 19
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(19, 1).Value = 0
$ws.Cells.Item(19, 2).Value = '강원삼척시FC이글스'
$ws.Cells.Item(19, 3).Value = '창단일자
20000101
주소
강원도 삼척시 동해대로 4234 (교동, 삼척종합운동장) 삼척시축구협회
연락처
010-6375-1414
감독
지중근
팀 관리에서 팀 SNS를 등록하세요'

# Row 20
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20, 1).Value = 0
$ws.Cells.Item(20, 2).Value = '강원삼척시가람FC'
$ws.Cells.Item(20, 3).Value = '창단일자
20200501
주소
강원 삼척시 중앙로4길 14-34 가온누리스포츠
연락처
010-7188-0104
감독
팀 관리에서 팀 SNS를 등록하세요'

# Row 21
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(21, 1).Value = 0
$ws.Cells.Item(21, 2).Value = '강원삼척시라온FC'
$ws.Cells.Item(21, 3).Value = '창단일자
20200129
주소
강원 삼척시 동해대로 4234 (교동, 삼척종합운동장)
연락처
010-6876-7720
감독
팀 관리에서 팀 SNS를 등록하세요'

# Row 22
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22, 1).Value = 0
$ws.Cells.Item(22, 2).Value = '강원삼척시미리내FC'
$ws.Cells.Item(22, 3).Value = '창단일자
20230102
주소
강원 삼척시 도계읍 장원길 63 (도계리, 경동아파트) 12동 307호
연락처
010-4614-1984
감독
김동진
팀 관리에서 팀 SNS를 등록하세요'

# Row 23
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(23, 1).Value = 0
$ws.Cells.Item(23, 2).Value = '강원삼척시삼척동자FC'
$ws.Cells.Item(23, 3).Value = '창단일자
20220711
주소
강원 삼척시 중앙시장1길 22-23 (남양동, 세연상가아파트) 708호
연락처
010-9528-2260
감독
팀 관리에서 팀 SNS를 등록하세요'

# Row 24
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A24").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(24, 1).Value = 0
$ws.Cells.Item(24, 2).Value = '강원삼척시미리내FC'
$ws.Cells.Item(24, 3).Value = '창단일자
20230102
주소
강원 삼척시 도계읍 장원길 63 (도계리, 경동아파트) 12동 307호
연락처
010-4614-1984
감독
김동진
팀 관리에서 팀 SNS를 등록하세요'

# Row 25
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(25, 1).Value = 0
$ws.Cells.Item(25, 2).Value = '강원삼척시원더플FC'
$ws.Cells.Item(25, 3).Value = '창단일자
20161002
주소
강원 삼척시 동해대로 4234 (교동, 삼척종합운동장) 삼척시축구협회
연락처
010-2758-1516
감독
이광철
팀 관리에서 팀 SNS를 등록하세요'

# Row 26
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A26").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(26, 1).Value = 0
$ws.Cells.Item(26, 2).Value = '강원속초시일오축구회'
$ws.Cells.Item(26, 3).Value = '창단일자
19810813
주소
강원도 속초시 번영로129번길 25 (동명동) (동명동)
연락처
010-4370-6085
감독
김명수
팀 관리에서 팀 SNS를 등록하세요'

# Row 27
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A27").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27, 1).Value = 0
$ws.Cells.Item(27, 2).Value = '강원속초시청해FC축구회'
$ws.Cells.Item(27, 3).Value = '창단일자
20180701
주소
강원 속초시 중앙로 96 (청학동) 로데오센터빌딩 4층
연락처
010-2803-1271
감독
박일균
팀 관리에서 팀 SNS를 등록하세요'

# Row 28
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A28").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(28, 1).Value = 0
$ws.Cells.Item(28, 2).Value = '강원속초시클래스'
$ws.Cells.Item(28, 3).Value = '창단일자
20160101
주소
강원 속초시 청대마을길 37 (조양동) 청대리막국수
연락처
010-3699-5210
감독
장유민
팀 관리에서 팀 SNS를 등록하세요'

# Row 29
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(29, 1).Value = 0
$ws.Cells.Item(29, 2).Value = '강원속초시한솔축구회'
$ws.Cells.Item(29, 3).Value = '창단일자
2001
주소
강원 속초시 청대마을길 16-1 (조양동)신태용 축구교실 속초센터
연락처
010-5389-9489
감독
이서호
팀 관리에서 팀 SNS를 등록하세요'

# Row 30
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(30, 1).Value = 0
$ws.Cells.Item(30, 2).Value = '강원속초시해마축구회'
$ws.Cells.Item(30, 3).Value = '창단일자
19990301
주소
강원도 속초시 청대로204번길 40 (조양동, 부영아파트5단지) 501동~508동 505동111호
연락처
010-3115-4285
감독
김병일
팀 관리에서 팀 SNS를 등록하세요'

# Row 31
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(31, 1).Value = 0
$ws.Cells.Item(31, 2).Value = '강원속초시현대축구회'
$ws.Cells.Item(31, 3).Value = '창단일자
19890101
주소
강원도 속초시 소평로 68 (조양동) 257-4
연락처
010-5362-7747
감독
김성백
팀 관리에서 팀 SNS를 등록하세요'

# Row 32
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(32, 1).Value = 0
$ws.Cells.Item(32, 2).Value = '강원양구군BAEKDUUTD'
$ws.Cells.Item(32, 3).Value = '창단일자
20211118
주소
강원 양구군 동면 금강산로 1801 (임당리, 대암아파트) e동 402호
연락처
010-2004-7995
감독
강민재
팀 관리에서 팀 SNS를 등록하세요'

# Row 33
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(33, 1).Value = 0
$ws.Cells.Item(33, 2).Value = '강원속초시현대축구회'
$ws.Cells.Item(33, 3).Value = '창단일자
19890101
주소
강원도 속초시 소평로 68 (조양동) 257-4
연락처
010-5362-7747
감독
김성백
팀 관리에서 팀 SNS를 등록하세요'

# Row 34
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(34, 1).Value = 0
$ws.Cells.Item(34, 2).Value = '강원영월군동강FC60대'
$ws.Cells.Item(34, 3).Value = '창단일자
20000301
주소
강원도 영월군 영월읍 제방안길 100 (하송리, 영월스포츠파크) 영월군체육회
연락처
010-5214-2290
감독
방극원1
팀 관리에서 팀 SNS를 등록하세요'
